# Update crypto price/volume table cells to match the latest scrape.
# Values in column D that look numeric are written with a leading
# apostrophe (Excel's quote-prefix) so they stay text cells instead of
# being auto-converted to numbers (matches the original inline-string cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.621.58"
$ws.Range("E2").Value = "  +3.34%  "
$ws.Range("D3").Value = "1.859.23"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'272.85"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.5278"
$ws.Range("E7").Value = "  +3.38%  "
$ws.Range("D8").Value = "'0.3378"
$ws.Range("D9").Value = "'0.06791"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").Value = "'19.83"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "'0.7922"
$ws.Range("E11").Value = "  -4.57%  "
$ws.Range("D12").Value = "'0.07732"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").Value = "1.831.06"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "'89.68"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "'5.121"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "'1.0000"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'14.40"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007990"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'0.9998"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "26.645.65"
$ws.Range("E20").Value = "  +3.25%  "
$ws.Range("D21").Value = "2.117.99"
$ws.Range("E21").Value = "  +3.66%  "
$ws.Range("D22").Value = "'4.714"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "'9.972"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "'6.093"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "'2.362"
$ws.Range("E25").Value = "  +5.19%  "
$ws.Range("D26").Value = "'145.86"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'112.00"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").Value = "'4.309"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "'4.301"
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").Value = "'0.08888"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").Value = "'0.04905"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'1.156"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").Value = "'0.7274"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "'2.881"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'3.224"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("D38").Value = "'2.326"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").Value = "'0.01843"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "'0.5065"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").Value = "'0.9402"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").Value = "'116.11"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "'6.132"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").Value = "'8.006"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "'0.9997"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'0.4397"
$ws.Range("E46").Value = "  -3.46%  "
$ws.Range("D47").Value = "'0.1324"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").Value = "'9.325"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'36.04"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").Value = "'0.05930"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").Value = "'1.470"
$ws.Range("E51").Value = "  -2.12%  "
